# #CRM-1615 Add Remarks in SF pending booking page
# Adds a new "Remarks" column (column Q) to the SF Pending Bookings list:
#   Q1 -> header "Remarks"                       (same look as the other
#                                                  bold/centered headers, e.g. P1)
#   Q2 -> placeholder "{bookings:booking_remarks}" (same look as the other
#                                                  left-aligned data cells, e.g. H2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell: Q1 = "Remarks" -------------------------------------
$ws.Range("Q1").Value = "Remarks"
# Match the formatting of the existing header row (bold, centered - same as P1)
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats

# --- New data cell: Q2 = "{bookings:booking_remarks}" --------------------
$ws.Range("Q2").Value = "{bookings:booking_remarks}"
# Match the formatting of the existing left-aligned data cells (same as H2)
$ws.Range("H2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# --- Restore the view/selection state left behind in the saved sheet -----
$null = $ws.Range("P4").Select()
